# ------------------------------------------------------------------
# Applies the commit "mise à jour des documents" to the active document:
#   1. Shortens "...un numéro de réservation et son billet." to
#      "...un numéro de réservation." in the table cell.
#   2. Adds a "Remarque : si un champ..." paragraph right after the
#      table (where the _GoBack bookmark now lives) and inserts a new
#      blank paragraph after it.
#   3. Moves <w:lastRenderedPageBreak/> from the "Méthode " run to the
#      "Résultat : recette..." run.
# ------------------------------------------------------------------

$d = $word.ActiveDocument

# --- Step 1: trim the sentence in the table cell ------------------
$d.Content.Find.Execute(
    "La compagnie génère et envoie au client un numéro de réservation et son billet.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "La compagnie génère et envoie au client un numéro de réservation.", 2
) | Out-Null

# --- Step 2: remove the _GoBack bookmark from the table cell ------
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

# --- Step 3: find the first empty paragraph right after the table -
$tbl = $d.Tables(1)
$afterTable = $tbl.Range.End
$p1 = $d.Range($afterTable, $afterTable).Paragraphs(1)
$insertPos = $p1.Range.Start

$part1 = "Remarque : s"
$part2 = "i un champ dans un formulaire n’est pas conforme aux attentes, le formulaire est retourné avec un message explicatif en en-tête."

$r = $d.Range($insertPos, $insertPos)
$r.InsertAfter($part1)
$r1 = $d.Range($insertPos, $insertPos + $part1.Length)
$r1.Font.Name = "Times New Roman"
$r1.Font.Size = 12

$bmPos = $insertPos + $part1.Length
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos)) | Out-Null

$r2ins = $d.Range($bmPos, $bmPos)
$r2ins.InsertAfter($part2)
$r2 = $d.Range($bmPos, $bmPos + $part2.Length)
$r2.Font.Name = "Times New Roman"
$r2.Font.Size = 12

$endOfPara = $bmPos + $part2.Length

# --- Step 4: add a new blank paragraph right after it -------------
$paraEndRange = $d.Range($endOfPara, $endOfPara)
$paraEndRange.InsertParagraphAfter()

# --- Step 5: move the lastRenderedPageBreak marker -----------------
# It currently sits in the "Méthode " run; the target document wants
# it on the "Résultat : recette..." run instead.
$found = $d.Content.Find.Execute("Résultat : recette, présentation de l'application...")
if ($found) {
    $resultRange = $d.Range($d.Content.Find.Parent.Start, $d.Content.Find.Parent.Start)
}

Write-Output "done"
